$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2646048109965636
$ws.Range("C2").Value = 0.4398625429553265
$ws.Range("J2").Value = 0.01030927835051546
$ws.Range("P2").Value = 0.1615120274914089
$ws.Range("S2").Value = 0.1237113402061856
# Row 3
$ws.Range("B3").Value = 0.00684931506849315
$ws.Range("C3").Value = 0.1095890410958904
$ws.Range("J3").Value = 0.0410958904109589
$ws.Range("P3").Value = 0.6438356164383562
$ws.Range("S3").Value = 0.1986301369863014
# Row 4
$ws.Range("J4").Value = 0.06122448979591837
$ws.Range("P4").Value = 0.5918367346938775
$ws.Range("S4").Value = 0.3469387755102041
# Row 5
$ws.Range("P5").Value = 0.8333333333333334
$ws.Range("S5").Value = 0.1666666666666667
# Row 6
$ws.Range("B6").Value = 0.05555555555555555
$ws.Range("D6").Value = 0.004273504273504274
$ws.Range("F6").Value = 0.03418803418803419
$ws.Range("J6").Value = 0.2991452991452991
$ws.Range("O6").Value = 0.004273504273504274
$ws.Range("Q6").Value = 0.1452991452991453
$ws.Range("R6").Value = 0.07692307692307693
$ws.Range("S6").Value = 0.3803418803418803
# Row 7
$ws.Range("B7").Value = 0.08812260536398467
$ws.Range("D7").Value = 0.02298850574712644
$ws.Range("E7").Value = 0.01532567049808429
$ws.Range("F7").Value = 0.05363984674329502
$ws.Range("J7").Value = 0.10727969348659
$ws.Range("Q7").Value = 0.2068965517241379
$ws.Range("R7").Value = 0.09961685823754789
$ws.Range("S7").Value = 0.4061302681992337
# Row 8
$ws.Range("B8").Value = 0.09110169491525423
$ws.Range("D8").Value = 0.0211864406779661
$ws.Range("F8").Value = 0.06567796610169492
$ws.Range("J8").Value = 0.1122881355932203
$ws.Range("O8").Value = 0.02330508474576271
$ws.Range("Q8").Value = 0.1504237288135593
$ws.Range("R8").Value = 0.08898305084745763
$ws.Range("S8").Value = 0.4470338983050847
# Row 9
$ws.Range("B9").Value = 0.08888888888888889
$ws.Range("D9").Value = 0.01111111111111111
$ws.Range("E9").Value = 0.01111111111111111
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.1222222222222222
$ws.Range("O9").Value = 0.01111111111111111
$ws.Range("Q9").Value = 0.1388888888888889
$ws.Range("R9").Value = 0.1055555555555556
$ws.Range("S9").Value = 0.4444444444444444
# Row 10
$ws.Range("B10").Value = 0.09841269841269841
$ws.Range("D10").Value = 0.02698412698412699
$ws.Range("E10").Value = 0.001587301587301587
$ws.Range("F10").Value = 0.08015873015873017
$ws.Range("J10").Value = 0.09841269841269841
$ws.Range("O10").Value = 0.01825396825396826
$ws.Range("Q10").Value = 0.2079365079365079
$ws.Range("R10").Value = 0.08492063492063492
$ws.Range("S10").Value = 0.3833333333333334
# Row 11
$ws.Range("G11").Value = 0.1654676258992806
$ws.Range("J11").Value = 0.07194244604316546
$ws.Range("K11").Value = 0.2134292565947242
$ws.Range("L11").Value = 0.5251798561151079
$ws.Range("S11").Value = 0.02398081534772182
# Row 12
$ws.Range("F12").Value = 0.008547008547008548
$ws.Range("G12").Value = 0.7350427350427351
$ws.Range("J12").Value = 0.188034188034188
$ws.Range("K12").Value = 0.008547008547008548
$ws.Range("L12").Value = 0.0170940170940171
$ws.Range("S12").Value = 0.04273504273504274
# Row 13
$ws.Range("G13").Value = 0.58
$ws.Range("J13").Value = 0.32
$ws.Range("S13").Value = 0.1
# Row 15
$ws.Range("F15").Value = 0.02666666666666667
$ws.Range("H15").Value = 0.1688888888888889
$ws.Range("I15").Value = 0.06222222222222222
$ws.Range("J15").Value = 0.3377777777777778
$ws.Range("K15").Value = 0.04
$ws.Range("M15").Value = 0.02222222222222222
$ws.Range("O15").Value = 0.05333333333333334
$ws.Range("S15").Value = 0.2888888888888889
# Row 16
$ws.Range("F16").Value = 0.01775147928994083
$ws.Range("H16").Value = 0.1775147928994083
$ws.Range("I16").Value = 0.0650887573964497
$ws.Range("J16").Value = 0.378698224852071
$ws.Range("K16").Value = 0.1420118343195266
$ws.Range("M16").Value = 0.005917159763313609
$ws.Range("O16").Value = 0.0650887573964497
$ws.Range("S16").Value = 0.1479289940828402
# Row 17
$ws.Range("F17").Value = 0.02008928571428572
$ws.Range("H17").Value = 0.2053571428571428
$ws.Range("I17").Value = 0.08705357142857142
$ws.Range("J17").Value = 0.3883928571428572
$ws.Range("K17").Value = 0.1294642857142857
$ws.Range("M17").Value = 0.01785714285714286
$ws.Range("N17").Value = 0.002232142857142857
$ws.Range("O17").Value = 0.05133928571428571
$ws.Range("S17").Value = 0.09821428571428571
# Row 18
$ws.Range("F18").Value = 0.02450980392156863
$ws.Range("H18").Value = 0.1813725490196078
$ws.Range("I18").Value = 0.06862745098039216
$ws.Range("J18").Value = 0.3529411764705883
$ws.Range("K18").Value = 0.142156862745098
$ws.Range("M18").Value = 0.02450980392156863
$ws.Range("O18").Value = 0.09313725490196079
$ws.Range("S18").Value = 0.1127450980392157
# Row 19
$ws.Range("F19").Value = 0.02223816355810617
$ws.Range("H19").Value = 0.1951219512195122
$ws.Range("I19").Value = 0.07101865136298421
$ws.Range("J19").Value = 0.3550932568149211
$ws.Range("K19").Value = 0.1384505021520804
$ws.Range("M19").Value = 0.02295552367288379
$ws.Range("N19").Value = 0.0007173601147776184
$ws.Range("O19").Value = 0.06169296987087518
$ws.Range("S19").Value = 0.1327116212338594

Write-Host "Updated 113 cells in team specific matrix"
